$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "RUBRO " (single concept column) is being split into two columns:
# "RUBRO TEMPORAL" and "RUBRO PERMANENTE". Insert a new column before the
# existing RUBRO column (B) so everything from B onward shifts one column right.
$ws.Range("B1").EntireColumn.Insert()

# Re-write the header row (row 1) with the new, wider set of columns.
$ws.Range("A1").Value = "NIT"
$ws.Range("B1").Value = "RUBRO TEMPORAL"
$ws.Range("C1").Value = "RUBRO PERMANENTE"
$ws.Range("D1").Value = "CONCEPTO"
$ws.Range("E1").Value = "UNIDAD 2"
$ws.Range("F1").Value = "UNIDAD 8"
$ws.Range("G1").Value = "UNIDAD 9 "
$ws.Range("H1").Value = " TOTAL "

# Size the columns to fit their (new) header contents.
$ws.Range("A1").EntireColumn.ColumnWidth = 3.667
$ws.Range("B1").EntireColumn.ColumnWidth = 21.667
$ws.Range("C1").EntireColumn.ColumnWidth = 25.0
$ws.Range("D1").EntireColumn.ColumnWidth = 13.333
$ws.Range("E1").EntireColumn.ColumnWidth = 10.834
$ws.Range("F1").EntireColumn.ColumnWidth = 10.834
$ws.Range("G1").EntireColumn.ColumnWidth = 11.5
$ws.Range("H1").EntireColumn.ColumnWidth = 8.667

# Select the whole sheet (as left by the author after the edit).
$ws.Cells.Select()
